$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "56.840.95"
$ws.Cells.Item(2, 5).Value = "  +0.30%  "
$ws.Cells.Item(3, 4).Value = "2.415.10"
$ws.Cells.Item(3, 5).Value = "  -3.42%  "
$ws.Cells.Item(4, 5).Value = "  +0.00%  "
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "486.59"
$ws.Cells.Item(5, 5).Value = "  -1.51%  "
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "153.65"
$ws.Cells.Item(6, 5).Value = "  +0.57%  "
$ws.Cells.Item(7, 2).Value = "USDC"
$ws.Cells.Item(7, 3).Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = "0.997"
$ws.Cells.Item(7, 5).Value = "  +0.00%  "
$ws.Cells.Item(8, 2).Value = "XRP"
$ws.Cells.Item(8, 3).Value = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = "0.613"
$ws.Cells.Item(8, 5).Value = "  +18.61%  "
$ws.Cells.Item(9, 4).Value = "2.427.34"
$ws.Cells.Item(9, 5).Value = "  -3.31%  "
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = "0.0998"
$ws.Cells.Item(10, 5).Value = "  +0.88%  "
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = "5.63"
$ws.Cells.Item(11, 5).Value = "  -2.73%  "
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = "0.335"
$ws.Cells.Item(12, 5).Value = "  -0.82%  "
$ws.Cells.Item(13, 5).Value = "  +0.86%  "
$ws.Cells.Item(14, 4).Value = "2.836.58"
$ws.Cells.Item(14, 5).Value = "  -3.26%  "
$ws.Cells.Item(15, 4).Value = "57.024.17"
$ws.Cells.Item(15, 5).Value = "  +0.55%  "
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = "20.74"
$ws.Cells.Item(16, 5).Value = "  -2.80%  "
$ws.Cells.Item(17, 5).Value = "  -2.96%  "
$ws.Cells.Item(18, 4).Value = "2.427.67"
$ws.Cells.Item(18, 5).Value = "  -3.42%  "
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = "4.75"
$ws.Cells.Item(19, 5).Value = "  +3.91%  "
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = "325.25"
$ws.Cells.Item(20, 5).Value = "  +0.21%  "
$ws.Cells.Item(21, 5).Value = "  -3.38%  "
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = "0.999"
$ws.Cells.Item(22, 5).Value = "  +0.07%  "
$ws.Cells.Item(23, 5).Value = "  -0.76%  "
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = "57.88"
$ws.Cells.Item(24, 5).Value = "  -1.51%  "
$ws.Cells.Item(25, 5).Value = "  -0.55%  "
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = "1.00"
$ws.Cells.Item(26, 5).Value = "  +0.26%  "
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = "0.157"
$ws.Cells.Item(27, 5).Value = "  -4.13%  "
$ws.Cells.Item(28, 4).Value = "2.529.07"
$ws.Cells.Item(28, 5).Value = "  -2.81%  "
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = "7.24"
$ws.Cells.Item(29, 5).Value = "  -6.01%  "
$ws.Cells.Item(30, 4).Value = "0.0₃0781"
$ws.Cells.Item(30, 5).Value = "  -3.68%  "
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = "0.999"
$ws.Cells.Item(31, 5).Value = "  +0.13%  "
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = "18.63"
$ws.Cells.Item(32, 5).Value = "  +0.98%  "
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = "148.58"
$ws.Cells.Item(33, 5).Value = "  -1.49%  "
$ws.Cells.Item(34, 5).Value = "  -0.27%  "
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = "5.29"
$ws.Cells.Item(35, 5).Value = "  +0.94%  "
$ws.Cells.Item(36, 5).Value = "  -2.08%  "
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = "3.70"
$ws.Cells.Item(37, 5).Value = "  -2.70%  "
$ws.Cells.Item(38, 5).Value = "  -3.75%  "
$ws.Cells.Item(39, 5).Value = "  +9.81%  "
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = "34.13"
$ws.Cells.Item(40, 5).Value = "  +0.32%  "
$ws.Cells.Item(41, 5).Value = "  -1.53%  "
$ws.Cells.Item(42, 5).Value = "  -0.04%  "
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = "0.994"
$ws.Cells.Item(43, 5).Value = "  +0.00%  "
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = "0.592"
$ws.Cells.Item(44, 5).Value = "  -3.59%  "
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = "266.56"
$ws.Cells.Item(45, 5).Value = "  +0.08%  "
$ws.Cells.Item(46, 5).Value = "  -5.50%  "
$ws.Cells.Item(47, 2).Value = "RenderToken"
$ws.Cells.Item(47, 3).Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = "4.69"
$ws.Cells.Item(47, 5).Value = "  -4.46%  "
$ws.Cells.Item(48, 2).Value = "WhiteBITCoin"
$ws.Cells.Item(48, 3).Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = "10.19"
$ws.Cells.Item(48, 5).Value = "  -0.37%  "
$ws.Cells.Item(49, 2).Value = "VeChain"
$ws.Cells.Item(49, 3).Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = "0.0227"
$ws.Cells.Item(49, 5).Value = "  -1.69%  "
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = "17.48"
$ws.Cells.Item(50, 5).Value = "  -2.86%  "
$ws.Cells.Item(51, 4).Value = "1.858.80"
$ws.Cells.Item(51, 5).Value = "  -2.77%  "
